{"js": "// Update the date line and the 25 \"three-digit \u00d7 one-digit\" answer cells.\n// Each old value is unique in the document, so searching the whole body\n// for each literal old string and replacing the single hit is sufficient.\n\nconst pairs = [\n  [\"2024-05-08 Wednesday\", \"2024-05-09 Thursday\"],\n  [\"922\u00d77=6454\", \"962\u00d79=8658\"],\n  [\"664\u00d74=2656\", \"609\u00d73=1827\"],\n  [\"727\u00d78=5816\", \"544\u00d78=4352\"],\n  [\"665\u00d76=3990\", \"605\u00d79=5445\"],\n  [\"405\u00d79=3645\", \"994\u00d73=2982\"],\n  [\"323\u00d76=1938\", \"983\u00d76=5898\"],\n  [\"132\u00d75=660\", \"400\u00d76=2400\"],\n  [\"498\u00d74=1992\", \"859\u00d75=4295\"],\n  [\"983\u00d72=1966\", \"275\u00d79=2475\"],\n  [\"318\u00d73=954\", \"929\u00d79=8361\"],\n  [\"999\u00d74=3996\", \"594\u00d73=1782\"],\n  [\"231\u00d74=924\", \"217\u00d74=868\"],\n  [\"663\u00d73=1989\", \"863\u00d76=5178\"],\n  [\"357\u00d72=714\", \"340\u00d72=680\"],\n  [\"638\u00d74=2552\", \"145\u00d78=1160\"],\n  [\"135\u00d77=945\", \"453\u00d77=3171\"],\n  [\"319\u00d79=2871\", \"919\u00d75=4595\"],\n  [\"326\u00d78=2608\", \"132\u00d73=396\"],\n  [\"254\u00d77=1778\", \"215\u00d79=1935\"],\n  [\"647\u00d76=3882\", \"294\u00d73=882\"],\n  [\"572\u00d76=3432\", \"866\u00d73=2598\"],\n  [\"359\u00d79=3231\", \"542\u00d73=1626\"],\n  [\"556\u00d77=3892\", \"658\u00d78=5264\"],\n  [\"562\u00d72=1124\", \"321\u00d76=1926\"],\n  [\"333\u00d79=2997\", \"667\u00d79=6003\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 \"three-digit \u00d7 one-digit\" answer cells.\n# Each old value is unique in the document, so a sequence of literal\n# (non-wildcard) Find/Replace passes over the whole document body is\n# sufficient and safe (no accidental re-matches between pairs).\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-05-08 Wednesday\", \"2024-05-09 Thursday\"),\n    @(\"922\u00d77=6454\", \"962\u00d79=8658\"),\n    @(\"664\u00d74=2656\", \"609\u00d73=1827\"),\n    @(\"727\u00d78=5816\", \"544\u00d78=4352\"),\n    @(\"665\u00d76=3990\", \"605\u00d79=5445\"),\n    @(\"405\u00d79=3645\", \"994\u00d73=2982\"),\n    @(\"323\u00d76=1938\", \"983\u00d76=5898\"),\n    @(\"132\u00d75=660\", \"400\u00d76=2400\"),\n    @(\"498\u00d74=1992\", \"859\u00d75=4295\"),\n    @(\"983\u00d72=1966\", \"275\u00d79=2475\"),\n    @(\"318\u00d73=954\", \"929\u00d79=8361\"),\n    @(\"999\u00d74=3996\", \"594\u00d73=1782\"),\n    @(\"231\u00d74=924\", \"217\u00d74=868\"),\n    @(\"663\u00d73=1989\", \"863\u00d76=5178\"),\n    @(\"357\u00d72=714\", \"340\u00d72=680\"),\n    @(\"638\u00d74=2552\", \"145\u00d78=1160\"),\n    @(\"135\u00d77=945\", \"453\u00d77=3171\"),\n    @(\"319\u00d79=2871\", \"919\u00d75=4595\"),\n    @(\"326\u00d78=2608\", \"132\u00d73=396\"),\n    @(\"254\u00d77=1778\", \"215\u00d79=1935\"),\n    @(\"647\u00d76=3882\", \"294\u00d73=882\"),\n    @(\"572\u00d76=3432\", \"866\u00d73=2598\"),\n    @(\"359\u00d79=3231\", \"542\u00d73=1626\"),\n    @(\"556\u00d77=3892\", \"658\u00d78=5264\"),\n    @(\"562\u00d72=1124\", \"321\u00d76=1926\"),\n    @(\"333\u00d79=2997\", \"667\u00d79=6003\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
